# Input binding localization: append new shared-string rows (85-173) to Sheet1
# describing mouse / gamepad / keyboard binding display names, and widen the
# existing red-highlight conditional formatting to the full column (Excel row max).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(85, 1).Value = 'LeftMouseButton'
$ws.Cells.Item(85, 2).Value = 'Left Mouse Button'
$ws.Rows.Item(85).RowHeight = 13.4

$ws.Cells.Item(86, 1).Value = 'MiddleMouseButton'
$ws.Cells.Item(86, 2).Value = 'Middle Mouse Button'
$ws.Rows.Item(86).RowHeight = 13.4

$ws.Cells.Item(87, 1).Value = 'RightMouseButton'
$ws.Cells.Item(87, 2).Value = 'Right Mouse Button'
$ws.Rows.Item(87).RowHeight = 13.4

$ws.Cells.Item(88, 1).Value = 'DPadUp'
$ws.Cells.Item(88, 2).Value = 'D-Pad Up'
$ws.Rows.Item(88).RowHeight = 13.4

$ws.Cells.Item(89, 1).Value = 'DPadDown'
$ws.Cells.Item(89, 2).Value = 'D-Pad Down'
$ws.Rows.Item(89).RowHeight = 13.4

$ws.Cells.Item(90, 1).Value = 'DPadLeft'
$ws.Cells.Item(90, 2).Value = 'D-Pad Left'
$ws.Rows.Item(90).RowHeight = 13.4

$ws.Cells.Item(91, 1).Value = 'DPadRight'
$ws.Cells.Item(91, 2).Value = 'D-Pad Right'
$ws.Rows.Item(91).RowHeight = 13.4

$ws.Cells.Item(92, 1).Value = 'Start'
$ws.Cells.Item(92, 2).Value = 'Start'
$ws.Rows.Item(92).RowHeight = 13.4

$ws.Cells.Item(93, 1).Value = 'Back'
$ws.Cells.Item(93, 2).Value = 'Back'
$ws.Rows.Item(93).RowHeight = 13.4

$ws.Cells.Item(94, 1).Value = 'LeftStick'
$ws.Cells.Item(94, 2).Value = 'Left Thumbstick Click'
$ws.Rows.Item(94).RowHeight = 13.4

$ws.Cells.Item(95, 1).Value = 'RightStick'
$ws.Cells.Item(95, 2).Value = 'Right Thumbstick Click'
$ws.Rows.Item(95).RowHeight = 13.4

$ws.Cells.Item(96, 1).Value = 'LeftShoulder'
$ws.Cells.Item(96, 2).Value = 'Left Shoulder'
$ws.Rows.Item(96).RowHeight = 13.4

$ws.Cells.Item(97, 1).Value = 'RightShoulder'
$ws.Cells.Item(97, 2).Value = 'Right Shoulder'
$ws.Rows.Item(97).RowHeight = 13.4

$ws.Cells.Item(98, 1).Value = 'LeftThumbstickLeft'
$ws.Cells.Item(98, 2).Value = 'Left Thumbstick Left'
$ws.Rows.Item(98).RowHeight = 13.4

$ws.Cells.Item(99, 1).Value = 'RightTrigger'
$ws.Cells.Item(99, 2).Value = 'Right Trigger'
$ws.Rows.Item(99).RowHeight = 13.4

$ws.Cells.Item(100, 1).Value = 'LeftTrigger'
$ws.Cells.Item(100, 2).Value = 'Left Trigger'
$ws.Rows.Item(100).RowHeight = 13.4

$ws.Cells.Item(101, 1).Value = 'RightThumbstickUp'
$ws.Cells.Item(101, 2).Value = 'Right Thumbstick Up'
$ws.Rows.Item(101).RowHeight = 13.4

$ws.Cells.Item(102, 1).Value = 'RightThumbstickDown'
$ws.Cells.Item(102, 2).Value = 'Right Thumbstick Down'
$ws.Rows.Item(102).RowHeight = 12.8
$ws.Cells.Item(102, 2).WrapText = $false

$ws.Cells.Item(103, 1).Value = 'RightThumbstickRight'
$ws.Cells.Item(103, 2).Value = 'Right Thumbstick Right'
$ws.Rows.Item(103).RowHeight = 12.8
$ws.Cells.Item(103, 2).WrapText = $false

$ws.Cells.Item(104, 1).Value = 'RightThumbstickLeft'
$ws.Cells.Item(104, 2).Value = 'Right Thumbstick Left'
$ws.Rows.Item(104).RowHeight = 12.8
$ws.Cells.Item(104, 2).WrapText = $false

$ws.Cells.Item(105, 1).Value = 'LeftThumbstickUp'
$ws.Cells.Item(105, 2).Value = 'Left Thumbstick Up'
$ws.Rows.Item(105).RowHeight = 12.8
$ws.Cells.Item(105, 2).WrapText = $false

$ws.Cells.Item(106, 1).Value = 'LeftThumbstickDown'
$ws.Cells.Item(106, 2).Value = 'Left Thumbstick Down'
$ws.Rows.Item(106).RowHeight = 12.8
$ws.Cells.Item(106, 2).WrapText = $false

$ws.Cells.Item(107, 1).Value = 'LeftThumbstickRight'
$ws.Cells.Item(107, 2).Value = 'Left Thumbstick Right'
$ws.Rows.Item(107).RowHeight = 12.8
$ws.Cells.Item(107, 2).WrapText = $false

$ws.Cells.Item(108, 1).Value = 'Tab'
$ws.Cells.Item(108, 2).Value = 'Tab'
$ws.Rows.Item(108).RowHeight = 13.4

$ws.Cells.Item(109, 1).Value = 'Enter'
$ws.Cells.Item(109, 2).Value = 'Enter'
$ws.Rows.Item(109).RowHeight = 13.4

$ws.Cells.Item(110, 1).Value = 'Pause'
$ws.Cells.Item(110, 2).Value = 'Pause'
$ws.Rows.Item(110).RowHeight = 13.4

$ws.Cells.Item(111, 1).Value = 'CapsLock'
$ws.Cells.Item(111, 2).Value = 'Capslock'
$ws.Rows.Item(111).RowHeight = 13.4

$ws.Cells.Item(112, 1).Value = 'Kana'
$ws.Cells.Item(112, 2).Value = 'Kana'
$ws.Rows.Item(112).RowHeight = 13.4

$ws.Cells.Item(113, 1).Value = 'Kanji'
$ws.Cells.Item(113, 2).Value = 'Kanji'
$ws.Rows.Item(113).RowHeight = 13.4

$ws.Cells.Item(114, 1).Value = 'Escape'
$ws.Cells.Item(114, 2).Value = 'Escape'
$ws.Rows.Item(114).RowHeight = 13.4

$ws.Cells.Item(115, 1).Value = 'Space'
$ws.Cells.Item(115, 2).Value = 'Spacebar'
$ws.Rows.Item(115).RowHeight = 13.4

$ws.Cells.Item(116, 1).Value = 'PageUp'
$ws.Cells.Item(116, 2).Value = 'Page Up'
$ws.Rows.Item(116).RowHeight = 13.4

$ws.Cells.Item(117, 1).Value = 'PageDown'
$ws.Cells.Item(117, 2).Value = 'Page Down'
$ws.Rows.Item(117).RowHeight = 13.4

$ws.Cells.Item(118, 1).Value = 'End'
$ws.Cells.Item(118, 2).Value = 'End'
$ws.Rows.Item(118).RowHeight = 13.4

$ws.Cells.Item(119, 1).Value = 'Home'
$ws.Cells.Item(119, 2).Value = 'Home'
$ws.Rows.Item(119).RowHeight = 13.4

$ws.Cells.Item(120, 1).Value = 'Left'
$ws.Cells.Item(120, 2).Value = 'Left'
$ws.Rows.Item(120).RowHeight = 13.4

$ws.Cells.Item(121, 1).Value = 'Up'
$ws.Cells.Item(121, 2).Value = 'Up'
$ws.Rows.Item(121).RowHeight = 13.4

$ws.Cells.Item(122, 1).Value = 'Right'
$ws.Cells.Item(122, 2).Value = 'Right'
$ws.Rows.Item(122).RowHeight = 13.4

$ws.Cells.Item(123, 1).Value = 'Down'
$ws.Cells.Item(123, 2).Value = 'Down'
$ws.Rows.Item(123).RowHeight = 13.4

$ws.Cells.Item(124, 1).Value = 'Insert'
$ws.Cells.Item(124, 2).Value = 'Insert'
$ws.Rows.Item(124).RowHeight = 13.4

$ws.Cells.Item(125, 1).Value = 'Delete'
$ws.Cells.Item(125, 2).Value = 'Delete'
$ws.Rows.Item(125).RowHeight = 13.4

$ws.Cells.Item(126, 1).Value = 'Help'
$ws.Cells.Item(126, 2).Value = 'Help'
$ws.Rows.Item(126).RowHeight = 13.4

$ws.Cells.Item(127, 1).Value = 'D0'
$ws.Cells.Item(127, 2).Value = 0
$ws.Rows.Item(127).RowHeight = 13.4

$ws.Cells.Item(128, 1).Value = 'D1'
$ws.Cells.Item(128, 2).Value = 1
$ws.Rows.Item(128).RowHeight = 13.4

$ws.Cells.Item(129, 1).Value = 'D2'
$ws.Cells.Item(129, 2).Value = 2
$ws.Rows.Item(129).RowHeight = 13.4

$ws.Cells.Item(130, 1).Value = 'D3'
$ws.Cells.Item(130, 2).Value = 3
$ws.Rows.Item(130).RowHeight = 13.4

$ws.Cells.Item(131, 1).Value = 'D4'
$ws.Cells.Item(131, 2).Value = 4
$ws.Rows.Item(131).RowHeight = 13.4

$ws.Cells.Item(132, 1).Value = 'D5'
$ws.Cells.Item(132, 2).Value = 5
$ws.Rows.Item(132).RowHeight = 13.4

$ws.Cells.Item(133, 1).Value = 'D6'
$ws.Cells.Item(133, 2).Value = 6
$ws.Rows.Item(133).RowHeight = 13.4

$ws.Cells.Item(134, 1).Value = 'D7'
$ws.Cells.Item(134, 2).Value = 7
$ws.Rows.Item(134).RowHeight = 13.4

$ws.Cells.Item(135, 1).Value = 'D8'
$ws.Cells.Item(135, 2).Value = 8
$ws.Rows.Item(135).RowHeight = 13.4

$ws.Cells.Item(136, 1).Value = 'D9'
$ws.Cells.Item(136, 2).Value = 9
$ws.Rows.Item(136).RowHeight = 13.4

$ws.Cells.Item(137, 1).Value = 'LeftWindows'
$ws.Cells.Item(137, 2).Value = 'Left Windows Key'
$ws.Rows.Item(137).RowHeight = 13.4

$ws.Cells.Item(138, 1).Value = 'RightWindows'
$ws.Cells.Item(138, 2).Value = 'Right Windows Key'
$ws.Rows.Item(138).RowHeight = 13.4

$ws.Cells.Item(139, 1).Value = 'NumPad0'
$ws.Cells.Item(139, 2).Value = 'Numpad 0'
$ws.Rows.Item(139).RowHeight = 13.4

$ws.Cells.Item(140, 1).Value = 'NumPad1'
$ws.Cells.Item(140, 2).Value = 'Numpad 1'
$ws.Rows.Item(140).RowHeight = 13.4

$ws.Cells.Item(141, 1).Value = 'NumPad2'
$ws.Cells.Item(141, 2).Value = 'Numpad 2'
$ws.Rows.Item(141).RowHeight = 13.4

$ws.Cells.Item(142, 1).Value = 'NumPad3'
$ws.Cells.Item(142, 2).Value = 'Numpad 3'
$ws.Rows.Item(142).RowHeight = 13.4

$ws.Cells.Item(143, 1).Value = 'NumPad4'
$ws.Cells.Item(143, 2).Value = 'Numpad 4'
$ws.Rows.Item(143).RowHeight = 13.4

$ws.Cells.Item(144, 1).Value = 'NumPad5'
$ws.Cells.Item(144, 2).Value = 'Numpad 5'
$ws.Rows.Item(144).RowHeight = 13.4

$ws.Cells.Item(145, 1).Value = 'NumPad6'
$ws.Cells.Item(145, 2).Value = 'Numpad 6'
$ws.Rows.Item(145).RowHeight = 13.4

$ws.Cells.Item(146, 1).Value = 'NumPad7'
$ws.Cells.Item(146, 2).Value = 'Numpad 7'
$ws.Rows.Item(146).RowHeight = 13.4

$ws.Cells.Item(147, 1).Value = 'NumPad8'
$ws.Cells.Item(147, 2).Value = 'Numpad 8'
$ws.Rows.Item(147).RowHeight = 13.4

$ws.Cells.Item(148, 1).Value = 'NumPad9'
$ws.Cells.Item(148, 2).Value = 'Numpad 9'
$ws.Rows.Item(148).RowHeight = 13.4

$ws.Cells.Item(149, 1).Value = 'Multiply'
$ws.Cells.Item(149, 2).Value = 'Numpad *'
$ws.Rows.Item(149).RowHeight = 13.4

$ws.Cells.Item(150, 1).Value = 'Add'
$ws.Cells.Item(150, 2).Value = 'Numpad +'
$ws.Rows.Item(150).RowHeight = 13.4

$ws.Cells.Item(151, 1).Value = 'Subtract'
$ws.Cells.Item(151, 2).Value = 'Numpad -'
$ws.Rows.Item(151).RowHeight = 13.4

$ws.Cells.Item(152, 1).Value = 'Decimal'
$ws.Cells.Item(152, 2).Value = 'Numpad .'
$ws.Rows.Item(152).RowHeight = 13.4

$ws.Cells.Item(153, 1).Value = 'Divide'
$ws.Cells.Item(153, 2).Value = 'Numpad /'
$ws.Rows.Item(153).RowHeight = 13.4

$ws.Cells.Item(154, 1).Value = 'OemClear'
$ws.Cells.Item(154, 2).Value = 'Clear'
$ws.Rows.Item(154).RowHeight = 12.8
$ws.Cells.Item(154, 2).WrapText = $false

$ws.Cells.Item(155, 1).Value = 'NumLock'
$ws.Cells.Item(155, 2).Value = 'NumLock'
$ws.Rows.Item(155).RowHeight = 13.4

$ws.Cells.Item(156, 1).Value = 'LeftShift'
$ws.Cells.Item(156, 2).Value = 'Left Shift'
$ws.Rows.Item(156).RowHeight = 13.4

$ws.Cells.Item(157, 1).Value = 'RightShift'
$ws.Cells.Item(157, 2).Value = 'Right Shift'
$ws.Rows.Item(157).RowHeight = 13.4

$ws.Cells.Item(158, 1).Value = 'LeftControl'
$ws.Cells.Item(158, 2).Value = 'Left Control'
$ws.Rows.Item(158).RowHeight = 13.4

$ws.Cells.Item(159, 1).Value = 'RightControl'
$ws.Cells.Item(159, 2).Value = 'Right Control'
$ws.Rows.Item(159).RowHeight = 13.4

$ws.Cells.Item(160, 1).Value = 'LeftAlt'
$ws.Cells.Item(160, 2).Value = 'Left Alt'
$ws.Rows.Item(160).RowHeight = 13.4

$ws.Cells.Item(161, 1).Value = 'RightAlt'
$ws.Cells.Item(161, 2).Value = 'Right Alt'
$ws.Rows.Item(161).RowHeight = 13.4

$ws.Cells.Item(162, 1).Value = 'OemBackslash'
$ws.Cells.Item(162, 2).Value = '\'
$ws.Rows.Item(162).RowHeight = 13.4

$ws.Cells.Item(163, 1).Value = 'OemSemicolon'
$ws.Cells.Item(163, 2).Value = ';'
$ws.Rows.Item(163).RowHeight = 13.4

$ws.Cells.Item(164, 1).Value = 'OemPipe'
$ws.Cells.Item(164, 2).Value = '\'
$ws.Rows.Item(164).RowHeight = 13.4

$ws.Cells.Item(165, 1).Value = 'OemPeriod'
$ws.Cells.Item(165, 2).Value = '.'
$ws.Rows.Item(165).RowHeight = 13.4

$ws.Cells.Item(166, 1).Value = 'OemComma'
$ws.Cells.Item(166, 2).Value = ','
$ws.Rows.Item(166).RowHeight = 13.4

$ws.Cells.Item(167, 1).Value = 'OemQuestion'
$ws.Cells.Item(167, 2).Value = '/'
$ws.Rows.Item(167).RowHeight = 13.4

$ws.Cells.Item(168, 1).Value = 'OemQuotes'
$ws.Cells.Item(168, 2).Value = ''''''
$ws.Rows.Item(168).RowHeight = 13.4

$ws.Cells.Item(169, 1).Value = 'OemOpenBrackets'
$ws.Cells.Item(169, 2).Value = '['
$ws.Rows.Item(169).RowHeight = 13.4

$ws.Cells.Item(170, 1).Value = 'OemCloseBrackets'
$ws.Cells.Item(170, 2).Value = ']'
$ws.Rows.Item(170).RowHeight = 13.4

$ws.Cells.Item(171, 1).Value = 'OemPlus'
$ws.Cells.Item(171, 2).Value = '''='
$ws.Rows.Item(171).RowHeight = 13.4

$ws.Cells.Item(172, 1).Value = 'OemMinus'
$ws.Cells.Item(172, 2).Value = '-'
$ws.Rows.Item(172).RowHeight = 13.4

$ws.Cells.Item(173, 1).Value = 'OemTilde'
$ws.Cells.Item(173, 2).Value = '`'
$ws.Rows.Item(173).RowHeight = 13.4

# The red "missing binding" highlight rule already covers C2:C1048575;
# extend it one row further so it spans the entire column, reusing the same rule.
$condRange = $ws.Range("C2:C1048575")
$fcs = $condRange.FormatConditions
if ($fcs.Count -ge 1) {
    $fcs.Item(1).ModifyAppliesToRange($ws.Range("C2:C1048576"))
}

